# Automatische test-sync: 2025-06-26 23:28:50
# Append a new testmail row (#9, "Ik heb een klacht") to the Logs sheet,
# bump the Dashboard's "Klacht / Probleem" category total, and extend the
# chart series / conditional-formatting ranges to cover the new rows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) Logs sheet: append row 41
# ---------------------------------------------------------------
$wsLogs = $wb.Worksheets.Item("Logs")

$newRow = 41
$wsLogs.Cells.Item($newRow, 1).Value = "Ik heb een klacht"
$wsLogs.Cells.Item($newRow, 2).Value = "mailmind.test@zohomail.eu"
$wsLogs.Cells.Item($newRow, 3).Value = "Testmail #9: Ik heb een klacht"
$wsLogs.Cells.Item($newRow, 4).Value = "Klacht / Probleem"
$wsLogs.Cells.Item($newRow, 5).Value = "Beste klant,`nBedankt voor het doorsturen van uw klacht. Om uw klacht zo goed mogelijk te kunnen behandelen, ontvangen wij graag meer informatie over de aard van de klacht. Kunt u ons meer details geven over wat er precies is misgegaan? `nMet vriendelijke groet,`n[Naam] `nKlantenservice Team"
$wsLogs.Cells.Item($newRow, 6).Value = "2025-06-26 23:28:44"
$wsLogs.Cells.Item($newRow, 7).Value = "Ja"
$wsLogs.Cells.Item($newRow, 8).Value = "Nee"
$wsLogs.Cells.Item($newRow, 9).Value = "Ja"

# Extend the conditional-formatting ranges (D/G/H/I) from row 40 to row 41
# so the newly appended row keeps getting highlighted like the others.
$condD = $wsLogs.Range("D2:D40").FormatConditions
for ($i = 1; $i -le $condD.Count; $i++) {
    $condD.Item($i).ModifyAppliesToRange($wsLogs.Range("D2:D41"))
}

$condG = $wsLogs.Range("G2:G40").FormatConditions
for ($i = 1; $i -le $condG.Count; $i++) {
    $condG.Item($i).ModifyAppliesToRange($wsLogs.Range("G2:G41"))
}

$condH = $wsLogs.Range("H2:H40").FormatConditions
for ($i = 1; $i -le $condH.Count; $i++) {
    $condH.Item($i).ModifyAppliesToRange($wsLogs.Range("H2:H41"))
}

$condI = $wsLogs.Range("I2:I40").FormatConditions
for ($i = 1; $i -le $condI.Count; $i++) {
    $condI.Item($i).ModifyAppliesToRange($wsLogs.Range("I2:I41"))
}

# ---------------------------------------------------------------
# 2) Dashboard sheet: add the "Klacht / Probleem" tally row
# ---------------------------------------------------------------
$wsDash = $wb.Worksheets.Item("Dashboard")
$wsDash.Cells.Item(8, 1).Value = "Klacht / Probleem"
$wsDash.Cells.Item(8, 2).Value = 1

# ---------------------------------------------------------------
# 3) Chart: extend the category/value series ranges to include row 8
# ---------------------------------------------------------------
$chart = $wsDash.ChartObjects(1).Chart
$series = $chart.SeriesCollection(1)
$series.XValues = "='Dashboard'!`$A`$2:`$A`$8"
$series.Values = "='Dashboard'!`$B`$2:`$B`$8"
